# Fixed errors for demo. Added eCommerceSelected intent
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New entity block: eCommerceName (rows 38-41) ------------------------
# Column A: entity name, merged down the block, centered horizontally.
$ws.Range("A38").Value = "eCommerceName"
$ws.Range("A38:A41").HorizontalAlignment = -4108   # xlCenter
$ws.Range("A38:A41").Merge()

# Column B: display values
$ws.Range("B38").Value = "Amazon"
$ws.Range("B39").Value = "MediaMarkt"
$ws.Range("B40").Value = "AliExpress"
$ws.Range("B41").Value = "GearBest"

# Column C: synonyms
$ws.Range("C38").Value = "amazon"
$ws.Range("C38").HorizontalAlignment = -4108
$ws.Range("C38").VerticalAlignment = -4108

$ws.Range("C39").Value = "mediamarkt"
$ws.Range("C39").HorizontalAlignment = -4108
$ws.Range("C39").VerticalAlignment = -4108

$ws.Range("C40").Value = "aliexpres`nali express`n"
$ws.Range("C40").WrapText = $true
$ws.Rows.Item(40).RowHeight = 45

$ws.Range("C41").Value = "gearbest`ngear best"
$ws.Range("C41").WrapText = $true
$ws.Rows.Item(41).RowHeight = 30

# --- View state: scroll position + selection ------------------------------
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C49:C50").Select()
